$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was reported for this market/product. It sits
# chronologically/logically at the top of this block of rows, so insert a
# fresh row at 69 (shifting the existing rows 69-81 down to 70-82) and
# populate it with the new observation.
$ws.Rows.Item(69).Insert()

$ws.Cells.Item(69, 1).Value = 5
$ws.Cells.Item(69, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(69, 3).Value = "Maule"
$ws.Cells.Item(69, 4).Value = 44889
$ws.Cells.Item(69, 5).Value = 7
$ws.Cells.Item(69, 6).Value = "Fruta"
$ws.Cells.Item(69, 7).Value = 100101
$ws.Cells.Item(69, 8).Value = "Berries"
$ws.Cells.Item(69, 9).Value = 100101001
$ws.Cells.Item(69, 10).Value = "Arándano (blue)"
$ws.Cells.Item(69, 11).Value = "Sin especificar"
$ws.Cells.Item(69, 12).Value = "Primera"
$ws.Cells.Item(69, 13).Value = 150
$ws.Cells.Item(69, 14).Value = 4000
$ws.Cells.Item(69, 15).Value = 4000
$ws.Cells.Item(69, 16).Value = 4000
$ws.Cells.Item(69, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(69, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(69, 19).Value = 2000
$ws.Cells.Item(69, 20).Value = 2
